# Auto-generated edit script: refreshes Chocobo_Profits market-price snapshot cells
# (currentAveragePrice[/NQ/HQ], LevePriceNQ/HQ, LeveProfitNQ/HQ) across all job sheets,
# matching the data pulled by the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 4079.762  # ALC H98: was 3600
$ws.Cells.Item(98, 9).Value = 1497.7273  # ALC I98: was 1080
$ws.Cells.Item(98, 10).Value = 6920  # ALC J98: was 5538.4614
$ws.Cells.Item(98, 11).Value = 1497.7273  # ALC K98: was 1080
$ws.Cells.Item(98, 12).Value = 6920  # ALC L98: was 5538.4614
$ws.Cells.Item(98, 13).Value = 0.2726999999999862  # ALC M98: was 418
$ws.Cells.Item(98, 14).Value = -9916  # ALC N98: was -8534.4614
$ws.Cells.Item(112, 8).Value = 1248  # ALC H112: was 1242.3871
$ws.Cells.Item(112, 10).Value = 1284.9656  # ALC J112: was 1278.4407
$ws.Cells.Item(112, 12).Value = 3854.8968  # ALC L112: was 3835.3221
$ws.Cells.Item(112, 14).Value = -6070.8968  # ALC N112: was -6051.3221
$ws.Cells.Item(114, 8).Value = 40000  # ALC H114: was 0
$ws.Cells.Item(114, 10).Value = 40000  # ALC J114: was 0
$ws.Cells.Item(114, 12).Value = 40000  # ALC L114: was 0
$ws.Cells.Item(114, 14).Value = -48678  # ALC N114: new value
$ws.Cells.Item(122, 8).Value = 4079.762  # ALC H122: was 3600
$ws.Cells.Item(122, 9).Value = 1497.7273  # ALC I122: was 1080
$ws.Cells.Item(122, 10).Value = 6920  # ALC J122: was 5538.4614
$ws.Cells.Item(122, 11).Value = 4493.1819  # ALC K122: was 3240
$ws.Cells.Item(122, 12).Value = 20760  # ALC L122: was 16615.3842
$ws.Cells.Item(122, 13).Value = -2043.1819  # ALC M122: was -790
$ws.Cells.Item(122, 14).Value = -25660  # ALC N122: was -21515.3842
$ws.Cells.Item(123, 8).Value = 41313.75  # ALC H123: was 41287.145
$ws.Cells.Item(123, 10).Value = 41751.668  # ALC J123: was 41802
$ws.Cells.Item(123, 12).Value = 41751.668  # ALC L123: was 41802
$ws.Cells.Item(123, 14).Value = -51551.668  # ALC N123: was -51602

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(6, 8).Value = 0  # ARM H6: was 19000
$ws.Cells.Item(6, 10).Value = 0  # ARM J6: was 19000
$ws.Cells.Item(6, 12).Value = 0  # ARM L6: was 19000
$ws.Cells.Item(6, 14).Value = $null  # ARM N6: clear (was -19346)
$ws.Cells.Item(114, 8).Value = 30398  # ARM H114: was 0
$ws.Cells.Item(114, 10).Value = 30398  # ARM J114: was 0
$ws.Cells.Item(114, 12).Value = 30398  # ARM L114: was 0
$ws.Cells.Item(114, 14).Value = -39076  # ARM N114: new value

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(110, 8).Value = 0  # BSM H110: was 50000
$ws.Cells.Item(110, 10).Value = 0  # BSM J110: was 50000
$ws.Cells.Item(110, 12).Value = 0  # BSM L110: was 50000
$ws.Cells.Item(110, 14).Value = $null  # BSM N110: clear (was -58180)
$ws.Cells.Item(114, 8).Value = 40684  # BSM H114: was 30342
$ws.Cells.Item(114, 10).Value = 40684  # BSM J114: was 30342
$ws.Cells.Item(114, 12).Value = 40684  # BSM L114: was 30342
$ws.Cells.Item(114, 14).Value = -49362  # BSM N114: was -39020
$ws.Cells.Item(118, 8).Value = 0  # BSM H118: was 28890
$ws.Cells.Item(118, 10).Value = 0  # BSM J118: was 28890
$ws.Cells.Item(118, 12).Value = 0  # BSM L118: was 28890
$ws.Cells.Item(118, 14).Value = $null  # BSM N118: clear (was -32204)
$ws.Cells.Item(132, 8).Value = 47632.5  # BSM H132: was 50195
$ws.Cells.Item(132, 10).Value = 47632.5  # BSM J132: was 50195
$ws.Cells.Item(132, 12).Value = 47632.5  # BSM L132: was 50195
$ws.Cells.Item(132, 14).Value = -57752.5  # BSM N132: was -60315

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(13, 8).Value = 11998  # CRP H13: was 15300
$ws.Cells.Item(13, 10).Value = 11998  # CRP J13: was 15300
$ws.Cells.Item(13, 12).Value = 11998  # CRP L13: was 15300
$ws.Cells.Item(13, 14).Value = -12276  # CRP N13: was -15578
$ws.Cells.Item(106, 8).Value = 40000  # CRP H106: was 37500
$ws.Cells.Item(106, 10).Value = 40000  # CRP J106: was 37500
$ws.Cells.Item(106, 12).Value = 40000  # CRP L106: was 37500
$ws.Cells.Item(106, 14).Value = -42524  # CRP N106: was -40024
$ws.Cells.Item(109, 8).Value = 35000  # CRP H109: was 34799.25
$ws.Cells.Item(109, 10).Value = 35000  # CRP J109: was 34799.25
$ws.Cells.Item(109, 12).Value = 35000  # CRP L109: was 34799.25
$ws.Cells.Item(109, 14).Value = -37080  # CRP N109: was -36879.25
$ws.Cells.Item(127, 8).Value = 32775.715  # CRP H127: was 41864
$ws.Cells.Item(127, 9).Value = 10000  # CRP I127: was 0
$ws.Cells.Item(127, 10).Value = 41886  # CRP J127: was 41864
$ws.Cells.Item(127, 11).Value = 10000  # CRP K127: was 0
$ws.Cells.Item(127, 12).Value = 41886  # CRP L127: was 41864
$ws.Cells.Item(127, 13).Value = -5040  # CRP M127: new value
$ws.Cells.Item(127, 14).Value = -51806  # CRP N127: was -51784
$ws.Cells.Item(137, 8).Value = 41435  # CRP H137: was 41447.5
$ws.Cells.Item(137, 10).Value = 41435  # CRP J137: was 41447.5
$ws.Cells.Item(137, 12).Value = 41435  # CRP L137: was 41447.5
$ws.Cells.Item(137, 14).Value = -51635  # CRP N137: was -51647.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 115.083336  # CUL H2: was 58.9375
$ws.Cells.Item(2, 9).Value = 60  # CUL I2: was 36.42857
$ws.Cells.Item(2, 10).Value = 133.44444  # CUL J2: was 76.44444
$ws.Cells.Item(2, 11).Value = 360  # CUL K2: was 218.57142
$ws.Cells.Item(2, 12).Value = 800.6666399999999  # CUL L2: was 458.66664
$ws.Cells.Item(2, 13).Value = -247  # CUL M2: was -105.57142
$ws.Cells.Item(2, 14).Value = -1026.66664  # CUL N2: was -684.66664
$ws.Cells.Item(12, 8).Value = 97.10526  # CUL H12: was 81.47369
$ws.Cells.Item(12, 9).Value = 22.5  # CUL I12: was 7
$ws.Cells.Item(12, 10).Value = 131.53847  # CUL J12: was 108.07143
$ws.Cells.Item(12, 11).Value = 67.5  # CUL K12: was 21
$ws.Cells.Item(12, 12).Value = 394.61541  # CUL L12: was 324.21429
$ws.Cells.Item(12, 13).Value = 105.5  # CUL M12: was 152
$ws.Cells.Item(12, 14).Value = -740.61541  # CUL N12: was -670.21429
$ws.Cells.Item(23, 8).Value = 232.42857  # CUL H23: was 246.71428
$ws.Cells.Item(23, 9).Value = 300  # CUL I23: was 0
$ws.Cells.Item(23, 10).Value = 221.16667  # CUL J23: was 246.71428
$ws.Cells.Item(23, 11).Value = 900  # CUL K23: was 0
$ws.Cells.Item(23, 12).Value = 663.50001  # CUL L23: was 740.14284
$ws.Cells.Item(23, 13).Value = -665  # CUL M23: new value
$ws.Cells.Item(23, 14).Value = -1133.50001  # CUL N23: was -1210.14284
$ws.Cells.Item(38, 8).Value = 108  # CUL H38: was 98.304344
$ws.Cells.Item(38, 9).Value = 96.28570999999999  # CUL I38: was 83.8125
$ws.Cells.Item(38, 11).Value = 288.85713  # CUL K38: was 251.4375
$ws.Cells.Item(38, 13).Value = 58.14287000000002  # CUL M38: was 95.5625
$ws.Cells.Item(122, 8).Value = 3106.1304  # CUL H122: was 3139.6
$ws.Cells.Item(122, 10).Value = 3535.3157  # CUL J122: was 3587.6216
$ws.Cells.Item(122, 12).Value = 31817.8413  # CUL L122: was 32288.5944
$ws.Cells.Item(122, 14).Value = -36717.8413  # CUL N122: was -37188.5944
$ws.Cells.Item(131, 8).Value = 692.83673  # CUL H131: was 670.5816
$ws.Cells.Item(131, 9).Value = 221.94737  # CUL I131: was 216.17392
$ws.Cells.Item(131, 10).Value = 806.0886  # CUL J131: was 809.93335
$ws.Cells.Item(131, 11).Value = 665.84211  # CUL K131: was 648.5217600000001
$ws.Cells.Item(131, 12).Value = 2418.2658  # CUL L131: was 2429.80005
$ws.Cells.Item(131, 13).Value = 4374.15789  # CUL M131: was 4391.47824
$ws.Cells.Item(131, 14).Value = -12498.2658  # CUL N131: was -12509.80005

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(5, 8).Value = 13500  # GSM H5: was 0
$ws.Cells.Item(5, 10).Value = 13500  # GSM J5: was 0
$ws.Cells.Item(5, 12).Value = 13500  # GSM L5: was 0
$ws.Cells.Item(5, 14).Value = -13724  # GSM N5: new value
$ws.Cells.Item(101, 8).Value = 30000  # GSM H101: was 50000
$ws.Cells.Item(101, 10).Value = 30000  # GSM J101: was 50000
$ws.Cells.Item(101, 12).Value = 30000  # GSM L101: was 50000
$ws.Cells.Item(101, 14).Value = -36490  # GSM N101: was -56490
$ws.Cells.Item(122, 8).Value = 4568.1333  # GSM H122: was 3500.7727
$ws.Cells.Item(122, 9).Value = 1771.6  # GSM I122: was 1587.5333
$ws.Cells.Item(122, 10).Value = 10161.2  # GSM J122: was 7600.5713
$ws.Cells.Item(122, 11).Value = 5314.799999999999  # GSM K122: was 4762.5999
$ws.Cells.Item(122, 12).Value = 30483.6  # GSM L122: was 22801.7139
$ws.Cells.Item(122, 13).Value = -2864.799999999999  # GSM M122: was -2312.5999
$ws.Cells.Item(122, 14).Value = -35383.60000000001  # GSM N122: was -27701.7139
$ws.Cells.Item(132, 8).Value = 2757.7827  # GSM H132: was 2776.2856
$ws.Cells.Item(132, 9).Value = 1960.7028  # GSM I132: was 2043.15
$ws.Cells.Item(132, 11).Value = 5882.1084  # GSM K132: was 6129.450000000001
$ws.Cells.Item(132, 13).Value = -3352.1084  # GSM M132: was -3599.450000000001
$ws.Cells.Item(134, 8).Value = 45758.223  # GSM H134: was 38853.617
$ws.Cells.Item(134, 10).Value = 45758.223  # GSM J134: was 38853.617
$ws.Cells.Item(134, 12).Value = 137274.669  # GSM L134: was 116560.851
$ws.Cells.Item(134, 14).Value = -142344.669  # GSM N134: was -121630.851

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(10, 8).Value = 31099.666  # LTW H10: was 31067.334
$ws.Cells.Item(10, 9).Value = 300  # LTW I10: was 203
$ws.Cells.Item(10, 11).Value = 300  # LTW K10: was 203
$ws.Cells.Item(10, 13).Value = -160  # LTW M10: was -63
$ws.Cells.Item(21, 8).Value = 47742.6  # LTW H21: was 20000
$ws.Cells.Item(21, 10).Value = 47742.6  # LTW J21: was 20000
$ws.Cells.Item(21, 12).Value = 47742.6  # LTW L21: was 20000
$ws.Cells.Item(21, 14).Value = -48090.6  # LTW N21: was -20348
$ws.Cells.Item(40, 8).Value = 4871.875  # LTW H40: was 5193.913
$ws.Cells.Item(40, 9).Value = 3888.4  # LTW I40: was 4388.923
$ws.Cells.Item(40, 10).Value = 6511  # LTW J40: was 6240.4
$ws.Cells.Item(40, 11).Value = 3888.4  # LTW K40: was 4388.923
$ws.Cells.Item(40, 12).Value = 6511  # LTW L40: was 6240.4
$ws.Cells.Item(40, 13).Value = -3752.4  # LTW M40: was -4252.923
$ws.Cells.Item(40, 14).Value = -6783  # LTW N40: was -6512.4
$ws.Cells.Item(46, 8).Value = 1579.3478  # LTW H46: was 1988.174
$ws.Cells.Item(46, 9).Value = 1042.8572  # LTW I46: was 1620
$ws.Cells.Item(46, 10).Value = 2413.889  # LTW J46: was 2678.5
$ws.Cells.Item(46, 11).Value = 1042.8572  # LTW K46: was 1620
$ws.Cells.Item(46, 12).Value = 2413.889  # LTW L46: was 2678.5
$ws.Cells.Item(46, 13).Value = -854.8571999999999  # LTW M46: was -1432
$ws.Cells.Item(46, 14).Value = -2789.889  # LTW N46: was -3054.5
$ws.Cells.Item(133, 8).Value = 30312.5  # LTW H133: was 30318.334
$ws.Cells.Item(133, 10).Value = 30312.5  # LTW J133: was 30318.334
$ws.Cells.Item(133, 12).Value = 30312.5  # LTW L133: was 30318.334
$ws.Cells.Item(133, 14).Value = -35372.5  # LTW N133: was -35378.334
$ws.Cells.Item(135, 8).Value = 0  # LTW H135: was 100000
$ws.Cells.Item(135, 10).Value = 0  # LTW J135: was 100000
$ws.Cells.Item(135, 12).Value = 0  # LTW L135: was 100000
$ws.Cells.Item(135, 14).Value = $null  # LTW N135: clear (was -110140)

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(46, 8).Value = 63401  # WVR H46: was 66326.125
$ws.Cells.Item(46, 10).Value = 63401  # WVR J46: was 66326.125
$ws.Cells.Item(46, 12).Value = 63401  # WVR L46: was 66326.125
$ws.Cells.Item(46, 14).Value = -63863  # WVR N46: was -66788.125
$ws.Cells.Item(98, 8).Value = 12795  # WVR H98: was 0
$ws.Cells.Item(98, 10).Value = 12795  # WVR J98: was 0
$ws.Cells.Item(98, 12).Value = 12795  # WVR L98: was 0
$ws.Cells.Item(98, 14).Value = -18785  # WVR N98: new value
$ws.Cells.Item(107, 8).Value = 548.73334  # WVR H107: was 603.8823
$ws.Cells.Item(107, 9).Value = 503.16666  # WVR I107: was 567.1539
$ws.Cells.Item(107, 10).Value = 731  # WVR J107: was 723.25
$ws.Cells.Item(107, 11).Value = 1509.49998  # WVR K107: was 1701.4617
$ws.Cells.Item(107, 12).Value = 2193  # WVR L107: was 2169.75
$ws.Cells.Item(107, 13).Value = 410.5000199999999  # WVR M107: was 218.5382999999999
$ws.Cells.Item(107, 14).Value = -6033  # WVR N107: was -6009.75
$ws.Cells.Item(122, 8).Value = 8949.9  # WVR H122: was 11249.75
$ws.Cells.Item(122, 9).Value = 6566.6665  # WVR I122: was 0
$ws.Cells.Item(122, 10).Value = 9971.286  # WVR J122: was 11249.75
$ws.Cells.Item(122, 11).Value = 19699.9995  # WVR K122: was 0
$ws.Cells.Item(122, 12).Value = 29913.858  # WVR L122: was 33749.25
$ws.Cells.Item(122, 13).Value = -17249.9995  # WVR M122: new value
$ws.Cells.Item(122, 14).Value = -34813.858  # WVR N122: was -38649.25
$ws.Cells.Item(132, 8).Value = 7579824.5  # WVR H132: was 8551547
$ws.Cells.Item(132, 9).Value = 3872.3142  # WVR I132: was 4452.6
$ws.Cells.Item(132, 11).Value = 11616.9426  # WVR K132: was 13357.8
$ws.Cells.Item(132, 13).Value = -9086.942599999998  # WVR M132: was -10827.8
$ws.Cells.Item(134, 8).Value = 63401  # WVR H134: was 66326.125
$ws.Cells.Item(134, 10).Value = 63401  # WVR J134: was 66326.125
$ws.Cells.Item(134, 12).Value = 190203  # WVR L134: was 198978.375
$ws.Cells.Item(134, 14).Value = -195273  # WVR N134: was -204048.375
